# Update of league bases (rows re-synced after odds provider re-scrape),
# per commit "Atualizacao de bases das ligas, do dia: 20-06-2024 as 20:11".
#
# The rows below keep the same match id in column A (row index) but the
# match records themselves (id, teams, odds, ...) have been re-ordered:
# some rows' B:AD content moved to a different row in this pass.
#
# Mapping of B:AD payloads (source row -> destination row):
#   237 -> 235
#   235 -> 236
#   236 -> 237
#   239 -> 238
#   238 -> 239
#   312 -> 309
#   313 -> 310
#   310 -> 312
#   309 -> 313
# (row 311 is left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot the B:AD payload of every source row before writing anything ---
$row235 = $ws.Range("B235:AD235").Value2
$row236 = $ws.Range("B236:AD236").Value2
$row237 = $ws.Range("B237:AD237").Value2
$row238 = $ws.Range("B238:AD238").Value2
$row239 = $ws.Range("B239:AD239").Value2

$row309 = $ws.Range("B309:AD309").Value2
$row310 = $ws.Range("B310:AD310").Value2
$row312 = $ws.Range("B312:AD312").Value2
$row313 = $ws.Range("B313:AD313").Value2

# --- write the rotated/swapped payloads back out ---
$ws.Range("B235:AD235").Value = $row237
$ws.Range("B236:AD236").Value = $row235
$ws.Range("B237:AD237").Value = $row236

$ws.Range("B238:AD238").Value = $row239
$ws.Range("B239:AD239").Value = $row238

$ws.Range("B309:AD309").Value = $row312
$ws.Range("B310:AD310").Value = $row313
$ws.Range("B312:AD312").Value = $row310
$ws.Range("B313:AD313").Value = $row309

Write-Output "rows re-synced"
